# "launch file on completion of scrip[t]"
#
# Turns off the "dashboard_1" flag (column G) for the inflation-series
# rows (CPI, PCE, Core CPI, Core PCE) so they no longer launch/open on
# script completion.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G5").Value = $false
$ws.Range("G6").Value = $false
$ws.Range("G7").Value = $false
$ws.Range("G8").Value = $false

# Return the cursor to the sheet's default cell.
$ws.Range("A1").Select()
